# "login logout with axios"
# - Switch workbook calculation to Manual (xlCalculationManual).
# - On the "Installation" sheet, append five new npm/install steps below the
#   existing "npm install pinia" entry (rows 85, 88, 90, 92, 94), matching the
#   blank-row spacing pattern already used throughout the sheet.
# - Rows 85/90/92 reuse the "copied cell" look already present on row 78
#   (Arial Unicode MS, 10pt) - pull that exact formatting across via
#   copy/paste-special so no new style entries are introduced.
# - Move the view/selection down to the newly added area.

$wb = $excel.ActiveWorkbook

# Workbook now recalculates manually.
$excel.Calculation = -4135   # xlCalculationManual

$ws = $wb.Worksheets.Item("Installation")
$ws.Activate() | Out-Null

$xlPasteFormats = -4122

$ws.Range("B85").Value = "npm install vue3-cookies --save"
$ws.Range("B78").Copy() | Out-Null
$ws.Range("B85").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B88").Value = " npm uninstall epic-spinners"

$ws.Range("B90").Value = "npm install vue-hooks"
$ws.Range("B78").Copy() | Out-Null
$ws.Range("B90").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B92").Value = "npm install vue-jwt-decode"
$ws.Range("B78").Copy() | Out-Null
$ws.Range("B92").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B94").Value = "npm install --save vue-axios"

$excel.CutCopyMode = $false

# Move the window's selection near the newly typed rows.
$ws.Range("F92").Select() | Out-Null
